# Updated cryptocurrency market data (Price + Volume(1h) change%) and
# corrected two pairs of swapped rows, as captured by the source commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.079.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.378.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.500'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.24'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0789'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.35'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.82%  '
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.750.09'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.398.92'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.07%  '
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.060.58'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.31'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.92%  '
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.56'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.46%  '
$ws.Range('E28').Value = '  +0.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.24'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0743'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.47'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.32%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.49%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.105'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.57%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.31'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.70%  '
$ws.Range('B39').Value = 'EnergySwap'
$ws.Range('C39').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +13.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.17%  '
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  -31.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.952.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0281'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('E45').Value = '  +2.31%  '
$ws.Range('E46').Value = '  -0.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -12.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.613.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '52.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.80%  '
